$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1: same style as the other header cells (e.g. E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("F1").Style = $ws.Range("E1").Style

# Data cells F2:F9: time_taken values (plain, unstyled like column E data cells)
$times = @(
    "2021-10-05 10:50:22.890621",
    "2021-10-05 10:50:22.890631",
    "2021-10-05 10:50:22.890634",
    "2021-10-05 10:50:22.890637",
    "2021-10-05 10:50:22.890640",
    "2021-10-05 10:50:22.890642",
    "2021-10-05 10:50:22.890645",
    "2021-10-05 10:50:22.890647"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $times[$i]
}
